# Update the "取得日時" (acquired-at) timestamp column for the appended
# batch of rows (rows 2-15) on the active sheet ("ランサーズ") to reflect
# the new fetch time: 2026-01-28 02:03:55.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2026-01-28 02:03:55"

for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
